$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, F and G hold date-like text (e.g. "2024-04-18") that must stay
# as plain text, not get auto-converted into Excel date serials. Force the
# column to Text format first, write the values, then restore the default
# "Normal" style so the cells end up with no explicit style, same as the rest.
$dateTextRanges = @("B2:B21", "F2:F21", "G2:G21")
foreach ($r in $dateTextRanges) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Cells.Item(2, 1).Value = "KB"
$ws.Cells.Item(2, 2).Value = "2024-04-18"
$ws.Cells.Item(2, 3).Value = "제일엠앤에스"
$ws.Cells.Item(2, 4).Value = "KB"
$ws.Cells.Item(2, 5).Value = "KB"
$ws.Cells.Item(2, 6).Value = "2024-04-23"
$ws.Cells.Item(2, 7).Value = "2024-04-30"
$ws.Cells.Item(2, 8).Value = 52800
$ws.Cells.Item(2, 9).Value = 2400000
$ws.Cells.Item(2, 10).Value = 22000
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 100

$ws.Cells.Item(3, 1).Value = "KB"
$ws.Cells.Item(3, 2).Value = "2024-05-07"
$ws.Cells.Item(3, 3).Value = "KB제28호스팩"
$ws.Cells.Item(3, 4).Value = "KB"
$ws.Cells.Item(3, 5).Value = "KB"
$ws.Cells.Item(3, 6).Value = "2024-05-10"
$ws.Cells.Item(3, 7).Value = "2024-05-17"
$ws.Cells.Item(3, 8).Value = 10000
$ws.Cells.Item(3, 9).Value = 5000000
$ws.Cells.Item(3, 10).Value = 2000
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 100

$ws.Cells.Item(4, 1).Value = "KB"
$ws.Cells.Item(4, 2).Value = "2024-04-25"
$ws.Cells.Item(4, 3).Value = "HD현대마린솔루션"
$ws.Cells.Item(4, 4).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Cells.Item(4, 5).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Cells.Item(4, 6).Value = "2024-04-30"
$ws.Cells.Item(4, 7).Value = "2024-05-08"
$ws.Cells.Item(4, 8).Value = 215255.4
$ws.Cells.Item(4, 9).Value = 8900000
$ws.Cells.Item(4, 10).Value = 83400
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 29

$ws.Cells.Item(5, 1).Value = "KB"
$ws.Cells.Item(5, 2).Value = "2024-04-23"
$ws.Cells.Item(5, 3).Value = "민테크"
$ws.Cells.Item(5, 4).Value = "KB"
$ws.Cells.Item(5, 5).Value = "KB"
$ws.Cells.Item(5, 6).Value = "2024-04-26"
$ws.Cells.Item(5, 7).Value = "2024-05-03"
$ws.Cells.Item(5, 8).Value = 31500
$ws.Cells.Item(5, 9).Value = 3000000
$ws.Cells.Item(5, 10).Value = 10500
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 100

$ws.Cells.Item(6, 1).Value = "NH"
$ws.Cells.Item(6, 2).Value = "2024-05-07"
$ws.Cells.Item(6, 3).Value = "아이씨티케이"
$ws.Cells.Item(6, 4).Value = "NH"
$ws.Cells.Item(6, 5).Value = "NH"
$ws.Cells.Item(6, 6).Value = "2024-05-10"
$ws.Cells.Item(6, 7).Value = "2024-05-17"
$ws.Cells.Item(6, 8).Value = 39400
$ws.Cells.Item(6, 9).Value = 1970000
$ws.Cells.Item(6, 10).Value = 20000
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 100

$ws.Cells.Item(7, 1).Value = "SK"
$ws.Cells.Item(7, 2).Value = "2024-04-23"
$ws.Cells.Item(7, 3).Value = "SK증권제12호스팩"
$ws.Cells.Item(7, 4).Value = "SK"
$ws.Cells.Item(7, 5).Value = "SK"
$ws.Cells.Item(7, 6).Value = "2024-04-26"
$ws.Cells.Item(7, 7).Value = "2024-05-07"
$ws.Cells.Item(7, 8).Value = 6000
$ws.Cells.Item(7, 9).Value = 3000000
$ws.Cells.Item(7, 10).Value = 2000
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 100

$ws.Cells.Item(8, 1).Value = "대신"
$ws.Cells.Item(8, 2).Value = "2024-04-25"
$ws.Cells.Item(8, 3).Value = "HD현대마린솔루션"
$ws.Cells.Item(8, 4).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Cells.Item(8, 5).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Cells.Item(8, 6).Value = "2024-04-30"
$ws.Cells.Item(8, 7).Value = "2024-05-08"
$ws.Cells.Item(8, 8).Value = 18556.5
$ws.Cells.Item(8, 9).Value = 8900000
$ws.Cells.Item(8, 10).Value = 83400
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 2.5

$ws.Cells.Item(9, 1).Value = "미래"
$ws.Cells.Item(9, 2).Value = "2024-05-20"
$ws.Cells.Item(9, 3).Value = "미래에셋비전스팩4호"
$ws.Cells.Item(9, 4).Value = "미래"
$ws.Cells.Item(9, 5).Value = "미래"
$ws.Cells.Item(9, 6).Value = "2024-05-23"
$ws.Cells.Item(9, 7).Value = "2024-05-29"
$ws.Cells.Item(9, 8).Value = 13300
$ws.Cells.Item(9, 9).Value = 6650000
$ws.Cells.Item(9, 10).Value = 2000
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 100

$ws.Cells.Item(10, 1).Value = "삼성"
$ws.Cells.Item(10, 2).Value = "2024-04-25"
$ws.Cells.Item(10, 3).Value = "HD현대마린솔루션"
$ws.Cells.Item(10, 4).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Cells.Item(10, 5).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Cells.Item(10, 6).Value = "2024-04-30"
$ws.Cells.Item(10, 7).Value = "2024-05-08"
$ws.Cells.Item(10, 8).Value = 18556.5
$ws.Cells.Item(10, 9).Value = 8900000
$ws.Cells.Item(10, 10).Value = 83400
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 2.5

$ws.Cells.Item(11, 1).Value = "삼성"
$ws.Cells.Item(11, 2).Value = "2024-05-13"
$ws.Cells.Item(11, 3).Value = "노브랜드"
$ws.Cells.Item(11, 4).Value = "삼성"
$ws.Cells.Item(11, 5).Value = "삼성"
$ws.Cells.Item(11, 6).Value = "2024-05-17"
$ws.Cells.Item(11, 7).Value = "2024-05-23"
$ws.Cells.Item(11, 8).Value = 16800
$ws.Cells.Item(11, 9).Value = 1200000
$ws.Cells.Item(11, 10).Value = 14000
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 100

$ws.Cells.Item(12, 1).Value = "신한"
$ws.Cells.Item(12, 2).Value = "2024-04-11"
$ws.Cells.Item(12, 3).Value = "신한제13호스팩"
$ws.Cells.Item(12, 4).Value = "신한"
$ws.Cells.Item(12, 5).Value = "신한"
$ws.Cells.Item(12, 6).Value = "2024-04-15"
$ws.Cells.Item(12, 7).Value = "2024-04-22"
$ws.Cells.Item(12, 8).Value = 6000
$ws.Cells.Item(12, 9).Value = 3000000
$ws.Cells.Item(12, 10).Value = 2000
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 100

$ws.Cells.Item(13, 1).Value = "신한"
$ws.Cells.Item(13, 2).Value = "2024-04-25"
$ws.Cells.Item(13, 3).Value = "HD현대마린솔루션"
$ws.Cells.Item(13, 4).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Cells.Item(13, 5).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Cells.Item(13, 6).Value = "2024-04-30"
$ws.Cells.Item(13, 7).Value = "2024-05-08"
$ws.Cells.Item(13, 8).Value = 74226
$ws.Cells.Item(13, 9).Value = 8900000
$ws.Cells.Item(13, 10).Value = 83400
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 10

$ws.Cells.Item(14, 1).Value = "신한"
$ws.Cells.Item(14, 2).Value = "2024-04-02"
$ws.Cells.Item(14, 3).Value = "신한제12호스팩"
$ws.Cells.Item(14, 4).Value = "신한"
$ws.Cells.Item(14, 5).Value = "신한"
$ws.Cells.Item(14, 6).Value = "2024-04-05"
$ws.Cells.Item(14, 7).Value = "2024-04-15"
$ws.Cells.Item(14, 8).Value = 10000
$ws.Cells.Item(14, 9).Value = 5000000
$ws.Cells.Item(14, 10).Value = 2000
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 100

$ws.Cells.Item(15, 1).Value = "유비에스"
$ws.Cells.Item(15, 2).Value = "2024-04-25"
$ws.Cells.Item(15, 3).Value = "HD현대마린솔루션"
$ws.Cells.Item(15, 4).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Cells.Item(15, 5).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Cells.Item(15, 6).Value = "2024-04-30"
$ws.Cells.Item(15, 7).Value = "2024-05-08"
$ws.Cells.Item(15, 8).Value = 170719.8
$ws.Cells.Item(15, 9).Value = 8900000
$ws.Cells.Item(15, 10).Value = 83400
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 23

$ws.Cells.Item(16, 1).Value = "유안타"
$ws.Cells.Item(16, 2).Value = "2024-04-22"
$ws.Cells.Item(16, 3).Value = "유안타제16호스팩"
$ws.Cells.Item(16, 4).Value = "유안타"
$ws.Cells.Item(16, 5).Value = "유안타"
$ws.Cells.Item(16, 6).Value = "2024-04-25"
$ws.Cells.Item(16, 7).Value = "2024-05-02"
$ws.Cells.Item(16, 8).Value = 10300
$ws.Cells.Item(16, 9).Value = 5150000
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 100

$ws.Cells.Item(17, 1).Value = "제이피모간회사"
$ws.Cells.Item(17, 2).Value = "2024-04-25"
$ws.Cells.Item(17, 3).Value = "HD현대마린솔루션"
$ws.Cells.Item(17, 4).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Cells.Item(17, 5).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Cells.Item(17, 6).Value = "2024-04-30"
$ws.Cells.Item(17, 7).Value = "2024-05-08"
$ws.Cells.Item(17, 8).Value = 170719.8
$ws.Cells.Item(17, 9).Value = 8900000
$ws.Cells.Item(17, 10).Value = 83400
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 23

$ws.Cells.Item(18, 1).Value = "하나"
$ws.Cells.Item(18, 2).Value = "2024-04-25"
$ws.Cells.Item(18, 3).Value = "HD현대마린솔루션"
$ws.Cells.Item(18, 4).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Cells.Item(18, 5).Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성"
$ws.Cells.Item(18, 6).Value = "2024-04-30"
$ws.Cells.Item(18, 7).Value = "2024-05-08"
$ws.Cells.Item(18, 8).Value = 74226
$ws.Cells.Item(18, 9).Value = 8900000
$ws.Cells.Item(18, 10).Value = 83400
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 10

$ws.Cells.Item(19, 1).Value = "하나"
$ws.Cells.Item(19, 2).Value = "2024-04-15"
$ws.Cells.Item(19, 3).Value = "하나33호스팩"
$ws.Cells.Item(19, 4).Value = "하나"
$ws.Cells.Item(19, 5).Value = "하나"
$ws.Cells.Item(19, 6).Value = "2024-04-18"
$ws.Cells.Item(19, 7).Value = "2024-04-24"
$ws.Cells.Item(19, 8).Value = 7000
$ws.Cells.Item(19, 9).Value = 3500000
$ws.Cells.Item(19, 10).Value = 2000
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 100

$ws.Cells.Item(20, 1).Value = "한국"
$ws.Cells.Item(20, 2).Value = "2024-04-24"
$ws.Cells.Item(20, 3).Value = "코칩"
$ws.Cells.Item(20, 4).Value = "한국"
$ws.Cells.Item(20, 5).Value = "한국"
$ws.Cells.Item(20, 6).Value = "2024-04-29"
$ws.Cells.Item(20, 7).Value = "2024-05-07"
$ws.Cells.Item(20, 8).Value = 27000
$ws.Cells.Item(20, 9).Value = 1500000
$ws.Cells.Item(20, 10).Value = 18000
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 100

$ws.Cells.Item(21, 1).Value = "한국"
$ws.Cells.Item(21, 2).Value = "2024-04-22"
$ws.Cells.Item(21, 3).Value = "디앤디파마텍"
$ws.Cells.Item(21, 4).Value = "한국"
$ws.Cells.Item(21, 5).Value = "한국"
$ws.Cells.Item(21, 6).Value = "2024-04-25"
$ws.Cells.Item(21, 7).Value = "2024-05-02"
$ws.Cells.Item(21, 8).Value = 36300
$ws.Cells.Item(21, 9).Value = 1100000
$ws.Cells.Item(21, 10).Value = 33000
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 100

foreach ($r in $dateTextRanges) {
    $ws.Range($r).Style = "Normal"
}
